$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the subject/month title in C1 (appended "-September" to the month range) ---
$ws.Range("C1").Value = "Power System Analysis and Design  (4th Semester Electrical)    EE-501                MONTH: AUGUST -September"

# --- Fill in the newly recorded "05_WEEK" (column H) attendance marks for every student row ---
$ws.Cells.Item(4, 8).Value = 44843

$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(6, 8).Value = 3
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(8, 8).Value = 3
$ws.Cells.Item(9, 8).Value = 3
$ws.Cells.Item(10, 8).Value = 3
$ws.Cells.Item(11, 8).Value = 3
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(14, 8).Value = 3
$ws.Cells.Item(15, 8).Value = 3
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(18, 8).Value = 3
$ws.Cells.Item(19, 8).Value = 3
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(21, 8).Value = 3
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(30, 8).Value = 3
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(32, 8).Value = 3
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(34, 8).Value = 3
$ws.Cells.Item(35, 8).Value = 3
$ws.Cells.Item(36, 8).Value = 3
$ws.Cells.Item(37, 8).Value = 3
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(41, 8).Value = 0

# --- Extend the "Number of Present Students" summary formula (row 42) across the E:H columns ---
$ws.Range("E42").Formula = "=SUM(E5:E41)/3"
$ws.Range("F42").Formula = "=SUM(F5:F41)/3"
$ws.Range("G42").Formula = "=SUM(G5:G41)/3"
$ws.Range("H42").Formula = "=SUM(H5:H41)/3"

# The fill-right also carried D42's formatting into the (still empty) I42:J42 cells
$ws.Range("D42").Copy() | Out-Null
$ws.Range("I42:J42").PasteSpecial(-4122) | Out-Null

# --- Drop the stray, entirely blank trailing row 43 ---
$ws.Rows.Item(43).Delete() | Out-Null

# --- Leave the selection where the editor last left it ---
$ws.Range("U3").Select() | Out-Null
